# Rename the worksheet "ZW1" -> "Template" (this workbook is now used as a
# template for worksheets created from sub-controllers), update the cell
# that echoed the sheet name, move the active selection, and widen column A
# to fit the new label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZW1")

# The label cell in A3 mirrored the old sheet/tab name - update it to match.
$ws.Range("A3").Value = "Template"

# Rename the sheet itself.
$ws.Name = "Template"

# Column A now needs to fit the "Template" label with best-fit sizing.
$ws.Columns.Item(1).AutoFit()

# Move/record the active selection as captured in the saved file.
$ws.Range("C11").Select()
